$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(45018.99999999999, 10, 9.999999988355944, 10.00000001377432),
    @(45074.99999999999, 10, 9.999999987162473, 10.00000001257481),
    @(45081.99999999999, 10, 9.999999986603028, 10.00000001324257),
    @(45088.99999999999, 10, 9.999999986373803, 10.00000001485364),
    @(45095.99999999999, 10, 9.999999983537933, 10.00000001680784),
    @(45102.99999999999, 10, 9.999999979216195, 10.00000002243636),
    @(45109.99999999999, 10, 9.999999975608503, 10.00000003155747),
    @(45116.99999999999, 10, 9.999999964148826, 10.0000000421282),
    @(45123.99999999999, 10, 9.999999952984515, 10.00000005709266),
    @(45130.99999999999, 10, 9.999999940850717, 10.00000007222099)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Range("A$r").Value = $row[0]
    $wsForecast.Range("B$r").Value = $row[1]
    $wsForecast.Range("C$r").Value = $row[2]
    $wsForecast.Range("D$r").Value = $row[3]
    $r++
}

# --- Copy formatting from the existing sheets so the new sheet matches style ---
# Header style (bold, centered, bordered) from row 1 of "Weekly Quantity"
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Date-number-format style from column A data cells of "Weekly Quantity"
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# Restore the originally active sheet
$wsWeekly.Activate() | Out-Null
